$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the typo'd username/email on row 4: "sarvn" -> "saravanakumar"
#    (A4 = username, C4 = OU-style login, F4 = email address).
$ws.Range("A4").Value = "saravanakumar"
$ws.Range("C4").Value = "saravanakumar.r"
$ws.Range("F4").Value = "saravanakumar.r@saravana.com"

# 2. Remove the hyperlinks that were attached to E4/F4, keeping the ones
#    on E3/F3 untouched. This engine's Hyperlink.Delete()/collection
#    Delete() only operates on the *entire* worksheet Hyperlinks
#    collection at once, so: remember the formatting + target info for
#    the links we want to keep, clear everything, then recreate just
#    those two (restoring their original cell style, since Add()
#    reformats the target range as a fresh hyperlink style).
$e3Style = $ws.Range("E3").Style
$f3Style = $ws.Range("F3").Style

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:D3skt0p@123", [Type]::Missing, [Type]::Missing, "D3skt0p@123")
$ws.Range("E3").Style = $e3Style

$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:saravana.r@saravana.com", [Type]::Missing, [Type]::Missing, "saravana.r@saravana.com")
$ws.Range("F3").Style = $f3Style

# 3. Update the sheet selection to A4:G4 with A4 as the active cell.
$ws.Range("A4:G4").Select()
